$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "Förändrad" (changed) date column C for rows 2-27 advances from
# 2023-09-15 (serial 45184) to 2023-09-16 (serial 45185).
for ($row = 2; $row -le 27; $row++) {
    $ws.Cells.Item($row, 3).Value = 45185
}
